$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 6757276
$ws.Cells.Item(2, 8).Value = 3
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 2.6
$ws.Cells.Item(2, 13).Value = 2.25
$ws.Cells.Item(2, 14).Value = 2.6
$ws.Cells.Item(2, 15).Value = 3.5
$ws.Cells.Item(2, 16).Value = 2.25
$ws.Cells.Item(2, 17).Value = 0.25
$ws.Cells.Item(2, 18).Value = 1.75
$ws.Cells.Item(2, 19).Value = 2.05
$ws.Cells.Item(2, 20).Value = 3.25
$ws.Cells.Item(2, 21).Value = 1.775
$ws.Cells.Item(2, 22).Value = 2.025
$ws.Cells.Item(2, 23).Value = 1.6
$ws.Cells.Item(2, 26).Value = 0.75
$ws.Cells.Item(2, 28).Value = -0.5
$ws.Cells.Item(2, 29).Value = 0.5125
# Row 3
$ws.Cells.Item(3, 2).Value = 6760228
$ws.Cells.Item(3, 8).Value = 4
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 11).Value = 2.25
$ws.Cells.Item(3, 13).Value = 2.625
$ws.Cells.Item(3, 14).Value = 1.75
$ws.Cells.Item(3, 15).Value = 3.6
$ws.Cells.Item(3, 16).Value = 3.8
$ws.Cells.Item(3, 17).Value = -0.5
$ws.Cells.Item(3, 18).Value = 1.8
$ws.Cells.Item(3, 19).Value = 2
$ws.Cells.Item(3, 20).Value = 2.75
$ws.Cells.Item(3, 21).Value = 1.8
$ws.Cells.Item(3, 22).Value = 2
$ws.Cells.Item(3, 23).Value = 0.75
$ws.Cells.Item(3, 26).Value = 0.8
$ws.Cells.Item(3, 28).Value = 0.8
$ws.Cells.Item(3, 29).Value = -1
# Row 5
$ws.Cells.Item(5, 7).Value = 'BSC Rapid Chemnitz'
# Row 11
$ws.Cells.Item(11, 2).Value = 7035048
$ws.Cells.Item(11, 7).Value = 'TuRU Dsseldorf'
$ws.Cells.Item(11, 8).Value = 1
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 'H'
$ws.Cells.Item(11, 11).Value = 3.25
$ws.Cells.Item(11, 12).Value = 4
$ws.Cells.Item(11, 13).Value = 1.8
$ws.Cells.Item(11, 14).Value = 2.9
$ws.Cells.Item(11, 15).Value = 4
$ws.Cells.Item(11, 16).Value = 1.95
$ws.Cells.Item(11, 17).Value = 0.5
$ws.Cells.Item(11, 20).Value = 3
$ws.Cells.Item(11, 21).Value = 1.75
$ws.Cells.Item(11, 22).Value = 1.95
$ws.Cells.Item(11, 23).Value = 1.9
$ws.Cells.Item(11, 25).Value = -1
$ws.Cells.Item(11, 26).Value = 0.8
$ws.Cells.Item(11, 27).Value = -1
$ws.Cells.Item(11, 29).Value = 0.95
# Row 12
$ws.Cells.Item(12, 2).Value = 7035046
$ws.Cells.Item(12, 7).Value = 'FC Viersen'
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 2
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 3.6
$ws.Cells.Item(12, 13).Value = 3
$ws.Cells.Item(12, 14).Value = 2
$ws.Cells.Item(12, 16).Value = 3
$ws.Cells.Item(12, 18).Value = 1.8
$ws.Cells.Item(12, 19).Value = 2
$ws.Cells.Item(12, 20).Value = 2.75
$ws.Cells.Item(12, 21).Value = 1.8
$ws.Cells.Item(12, 22).Value = 2
$ws.Cells.Item(12, 25).Value = 2
$ws.Cells.Item(12, 27).Value = 1
$ws.Cells.Item(12, 28).Value = -1
$ws.Cells.Item(12, 29).Value = 1
# Row 13
$ws.Cells.Item(13, 2).Value = 7035047
$ws.Cells.Item(13, 7).Value = 'VfL Viktoria JuchenGarzweiler'
$ws.Cells.Item(13, 8).Value = 3
$ws.Cells.Item(13, 9).Value = 4
$ws.Cells.Item(13, 10).Value = 'A'
$ws.Cells.Item(13, 11).Value = 1.909
$ws.Cells.Item(13, 12).Value = 3.75
$ws.Cells.Item(13, 13).Value = 3.1
$ws.Cells.Item(13, 14).Value = 2.2
$ws.Cells.Item(13, 15).Value = 3.6
$ws.Cells.Item(13, 16).Value = 2.625
$ws.Cells.Item(13, 17).Value = -0.25
$ws.Cells.Item(13, 18).Value = 2
$ws.Cells.Item(13, 19).Value = 1.8
$ws.Cells.Item(13, 21).Value = 1.825
$ws.Cells.Item(13, 22).Value = 1.975
$ws.Cells.Item(13, 23).Value = -1
$ws.Cells.Item(13, 25).Value = 1.625
$ws.Cells.Item(13, 26).Value = -1
$ws.Cells.Item(13, 27).Value = 0.8
$ws.Cells.Item(13, 28).Value = 0.825
$ws.Cells.Item(13, 29).Value = -1
# Row 14
$ws.Cells.Item(14, 2).Value = 7089911
$ws.Cells.Item(14, 7).Value = 'VfB Frohnhausen'
$ws.Cells.Item(14, 8).Value = 4
$ws.Cells.Item(14, 9).Value = 3
$ws.Cells.Item(14, 11).Value = 2.25
$ws.Cells.Item(14, 12).Value = 3.75
$ws.Cells.Item(14, 13).Value = 2.5
$ws.Cells.Item(14, 14).Value = 2.25
$ws.Cells.Item(14, 15).Value = 3.75
$ws.Cells.Item(14, 16).Value = 2.5
$ws.Cells.Item(14, 17).Value = 0
$ws.Cells.Item(14, 18).Value = 1.8
$ws.Cells.Item(14, 19).Value = 2
$ws.Cells.Item(14, 20).Value = 3.5
$ws.Cells.Item(14, 21).Value = 1.8
$ws.Cells.Item(14, 22).Value = 2
$ws.Cells.Item(14, 23).Value = 1.25
$ws.Cells.Item(14, 26).Value = 0.8
$ws.Cells.Item(14, 28).Value = 0.8
$ws.Cells.Item(14, 29).Value = -1
# Row 15
$ws.Cells.Item(15, 2).Value = 7089910
$ws.Cells.Item(15, 7).Value = 'TuRU Dsseldorf'
$ws.Cells.Item(15, 8).Value = 2
$ws.Cells.Item(15, 9).Value = 1
$ws.Cells.Item(15, 11).Value = 3.25
$ws.Cells.Item(15, 12).Value = 4
$ws.Cells.Item(15, 13).Value = 1.8
$ws.Cells.Item(15, 14).Value = 3.25
$ws.Cells.Item(15, 15).Value = 4
$ws.Cells.Item(15, 16).Value = 1.8
$ws.Cells.Item(15, 17).Value = 0.5
$ws.Cells.Item(15, 18).Value = 1.975
$ws.Cells.Item(15, 19).Value = 1.825
$ws.Cells.Item(15, 20).Value = 3.25
$ws.Cells.Item(15, 21).Value = 1.85
$ws.Cells.Item(15, 22).Value = 1.95
$ws.Cells.Item(15, 23).Value = 2.25
$ws.Cells.Item(15, 26).Value = 0.9750000000000001
$ws.Cells.Item(15, 28).Value = -0.5
$ws.Cells.Item(15, 29).Value = 0.475
# Row 18
$ws.Cells.Item(18, 6).Value = 'Spvgg Steele 0309'
# Row 21
$ws.Cells.Item(21, 6).Value = 'SC Dsseldorf West'
# Row 25
$ws.Cells.Item(25, 7).Value = 'SC Dsseldorf West'
# Row 27
$ws.Cells.Item(27, 6).Value = 'Spvgg Steele 0309'
# Row 40
$ws.Cells.Item(40, 2).Value = 7465503
$ws.Cells.Item(40, 6).Value = 'FC Monheim'
$ws.Cells.Item(40, 7).Value = 'FC Remscheid'
$ws.Cells.Item(40, 8).Value = 2
$ws.Cells.Item(40, 11).Value = 1.909
$ws.Cells.Item(40, 13).Value = 3.1
$ws.Cells.Item(40, 14).Value = 1.615
$ws.Cells.Item(40, 15).Value = 4
$ws.Cells.Item(40, 16).Value = 4.2
$ws.Cells.Item(40, 17).Value = -0.75
$ws.Cells.Item(40, 18).Value = 1.825
$ws.Cells.Item(40, 19).Value = 1.975
$ws.Cells.Item(40, 20).Value = 3.25
$ws.Cells.Item(40, 21).Value = 1.95
$ws.Cells.Item(40, 22).Value = 1.85
$ws.Cells.Item(40, 23).Value = 0.615
$ws.Cells.Item(40, 26).Value = 0.825
$ws.Cells.Item(40, 28).Value = -1
$ws.Cells.Item(40, 29).Value = 0.8500000000000001
# Row 41
$ws.Cells.Item(41, 2).Value = 7465502
$ws.Cells.Item(41, 6).Value = 'FC Pesch 1956'
$ws.Cells.Item(41, 7).Value = 'SSV Bornheim'
$ws.Cells.Item(41, 8).Value = 4
$ws.Cells.Item(41, 11).Value = 2
$ws.Cells.Item(41, 13).Value = 2.875
$ws.Cells.Item(41, 14).Value = 2
$ws.Cells.Item(41, 15).Value = 3.75
$ws.Cells.Item(41, 16).Value = 2.9
$ws.Cells.Item(41, 17).Value = -0.25
$ws.Cells.Item(41, 18).Value = 1.8
$ws.Cells.Item(41, 19).Value = 2
$ws.Cells.Item(41, 20).Value = 3
$ws.Cells.Item(41, 21).Value = 1.775
$ws.Cells.Item(41, 22).Value = 2.025
$ws.Cells.Item(41, 23).Value = 1
$ws.Cells.Item(41, 26).Value = 0.8
$ws.Cells.Item(41, 28).Value = 0.7749999999999999
$ws.Cells.Item(41, 29).Value = -1
# Row 49
$ws.Cells.Item(49, 6).Value = 'Spvgg Steele 0309'
# Row 55
$ws.Cells.Item(55, 6).Value = 'Cronenberger SC'
# Row 58
$ws.Cells.Item(58, 7).Value = 'ASV Mettmann'
# Row 61
$ws.Cells.Item(61, 6).Value = 'BSC Rapid Chemnitz'
# Row 62
$ws.Cells.Item(62, 6).Value = 'SG Unterrath'
# Row 69
$ws.Cells.Item(69, 6).Value = 'SC Dsseldorf West'
# Row 71
$ws.Cells.Item(71, 7).Value = 'SC Dsseldorf West'
# Row 74
$ws.Cells.Item(74, 6).Value = 'ASV Mettmann'
